# Correction in SA algorithm and 746 logs
# Updates the Fitness values (column C) of the run_19 log sheet to reflect
# the corrected simulated-annealing run values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12 (Generation 0-10): Fitness corrected to 8234
$ws.Range("C2:C12").Value = 8234

# Rows 13-37 (Generation 11-35): Fitness corrected to 7657
$ws.Range("C13:C37").Value = 7657

# Rows 38-252 (Generation 36-250): Fitness corrected to 7573
$ws.Range("C38:C252").Value = 7573
